# "little adventure improvemnt corrected"
#
# Rows 2-5 of Sheet1 had their "location" (B) and "Temperature" (C) columns
# both pointing at the same garbled leftover Selenium/ChromeDriver debug
# string (e.g. "16.04.$[[ChromeDriver: chrome on XP (...)] -> id: xPat]").
# This restores the real coordinate pair in column B and the real
# temperature reading in column C for dusseldorf, Nice, Marseille and
# Monte Carlo.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# dusseldorf
$ws.Range("B2").Value = "51.2217,6.7762"
$ws.Range("C2").Value = "15.82."

# Nice
$ws.Range("B3").Value = "44,7.25"
$ws.Range("C3").Value = "13.34."

# Marseille
$ws.Range("B4").Value = "43.3333,5.5"
$ws.Range("C4").Value = "20.24."

# Monte Carlo
$ws.Range("B5").Value = "43.7496,7.437"
$ws.Range("C5").Value = "21.83."

# The corrected (shorter) values nudge the best-fit width of columns B/C
# down slightly from their old value.
$ws.Columns.Item(2).ColumnWidth = 81.67
$ws.Columns.Item(3).ColumnWidth = 81.67
